# Generate Report for Handback
# Replace the two UUID-named file references and refresh their associated
# timestamps across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$newA = "a64fdfab-c11e-4542-aa3d-e5683b07d293"
$newB = "ffffa11e9b76-c798-425b-984d-8f94529aae2c"

$newAXlfZh = "$newA.4c04c08c4c27ec13b3c2344d1a97a843ae0247b4.zh-cn.xlf"
$newAXlfDe = "$newA.4c04c08c4c27ec13b3c2344d1a97a843ae0247b4.de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newA.md"
$ws.Range("B2").Value = "e2e\$newA.md"
$ws.Range("G2").Value = "2016-08-30 21:11:15"

$ws.Range("A3").Value = "$newB.md"
$ws.Range("B3").Value = "e2e\$newB.md"
$ws.Range("G3").Value = "2016-08-30 21:11:15"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = "e2e\$newA.md"
    } elseif ($addr -eq '$B$3') {
        $h.TextToDisplay = "e2e\$newB.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newA.md"
$ws.Range("I2").Value = "$newA.md"
$ws.Range("G2").Value = $newAXlfZh
$ws.Range("H2").Value = "2016-08-30 21:11:00"
$ws.Range("J2").Value = $newAXlfZh
$ws.Range("K2").Value = "2016-08-30 21:11:32"

$ws.Range("A3").Value = "$newB.md"
$ws.Range("I3").Value = "$newB.md"
$ws.Range("G3").Value = $newAXlfZh
$ws.Range("H3").Value = "2016-08-30 21:11:00"
$ws.Range("J3").Value = $newAXlfZh
$ws.Range("K3").Value = "2016-08-30 21:11:32"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newA.md"
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = "$newA.md"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "$newB.md"
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = "$newB.md"
    }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newA.md"
$ws.Range("I2").Value = "$newA.md"
$ws.Range("G2").Value = $newAXlfDe
$ws.Range("H2").Value = "2016-08-30 21:11:15"
$ws.Range("J2").Value = $newAXlfDe
$ws.Range("K2").Value = "2016-08-30 21:11:40"

$ws.Range("A3").Value = "$newB.md"
$ws.Range("I3").Value = "$newB.md"
$ws.Range("G3").Value = $newAXlfDe
$ws.Range("H3").Value = "2016-08-30 21:11:15"
$ws.Range("J3").Value = $newAXlfDe
$ws.Range("K3").Value = "2016-08-30 21:11:40"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newA.md"
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = "$newA.md"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "$newB.md"
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = "$newB.md"
    }
}
